# Apply updated cryptocurrency price/volume data to worksheet cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.016.79"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.596.51"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.479"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.248"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("D12").Value = "1.818.56"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "1.592.33"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "26.007.28"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("E24").Value = "  +13.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E27").Value = "  -7.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("D36").Value = "1.129.43"
$ws.Range("E36").Value = "  +2.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.800"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.492"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").Value = "1.729.36"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "53.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.406"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "0.0₇0944"
$ws.Range("E51").Value = "  -15.53%  "
